# Re-order the "Recorded By" (column G) list of recorder names/emails on the
# "Session Analysis Results" sheet:
#   - if the comma-separated list contains the literal, case-sensitive
#     token "System", move it to the front (keeping the relative order of
#     the remaining entries), e.g.
#       "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#       "system, System, backup@backdoor.com" -> "System, system, backup@backdoor.com"
#   - otherwise, if there are exactly two entries, swap them, e.g.
#       "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"
#   - single-entry cells (and any other shape) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column 7 = G ("Recorded By"); row 1 is the header, data starts at row 2.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $parts = @($current -split ", ")

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newValue = [string]::Join(", ", $newParts)
    }
    elseif ($parts.Count -eq 2) {
        $newValue = [string]::Join(", ", @($parts[1], $parts[0]))
    }
    else {
        $newValue = $current
    }

    if (-not $newValue.Equals($current)) {
        $cell.Value = $newValue
    }
}
